# Update "想去人数" (F column) counts across sheets to reflect the latest
# scraped data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 249
$ws1.Range("F3").Value  = 442
$ws1.Range("F6").Value  = 539
$ws1.Range("F9").Value  = 266
$ws1.Range("F10").Value = 376
$ws1.Range("F12").Value = 625
$ws1.Range("F13").Value = 746
$ws1.Range("F14").Value = 1499
$ws1.Range("F15").Value = 1499
$ws1.Range("F17").Value = 26
$ws1.Range("F18").Value = 1345
$ws1.Range("F20").Value = 283
$ws1.Range("F23").Value = 96
$ws1.Range("F24").Value = 6532
$ws1.Range("F25").Value = 4862
$ws1.Range("F26").Value = 137
$ws1.Range("F29").Value = 149
$ws1.Range("F30").Value = 27
$ws1.Range("F32").Value = 1270
$ws1.Range("F33").Value = 188
$ws1.Range("F34").Value = 241
$ws1.Range("F35").Value = 594
$ws1.Range("F38").Value = 235
$ws1.Range("F39").Value = 148
$ws1.Range("F43").Value = 93

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 10

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 167
$ws3.Range("F3").Value = 2442
$ws3.Range("F4").Value = 188
$ws3.Range("F5").Value = 46

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 249
$ws4.Range("F3").Value  = 167
$ws4.Range("F4").Value  = 442
$ws4.Range("F7").Value  = 188
$ws4.Range("F8").Value  = 46
$ws4.Range("F10").Value = 539
$ws4.Range("F13").Value = 266
$ws4.Range("F15").Value = 376
$ws4.Range("F17").Value = 625
$ws4.Range("F18").Value = 746
$ws4.Range("F19").Value = 1499
$ws4.Range("F20").Value = 1499
$ws4.Range("F22").Value = 26
$ws4.Range("F23").Value = 1345
$ws4.Range("F25").Value = 283
$ws4.Range("F27").Value = 96
$ws4.Range("F30").Value = 6532
$ws4.Range("F31").Value = 4862
$ws4.Range("F32").Value = 137
$ws4.Range("F33").Value = 27
$ws4.Range("F34").Value = 1270
$ws4.Range("F35").Value = 188
$ws4.Range("F36").Value = 241
$ws4.Range("F38").Value = 594
$ws4.Range("F44").Value = 235
$ws4.Range("F46").Value = 59
$ws4.Range("F48").Value = 93
